# Add data for 2024-05-21
# Updates 2024 year-to-date (column K) crime totals across Citywide,
# By Neighborhood, and individual neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 2862
$ws.Range("K3").Value = 2801
$ws.Range("K4").Value = 581
$ws.Range("K5").Value = 181
$ws.Range("K6").Value = 3416
$ws.Range("K7").Value = 9841

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 76
$ws.Range("K7").Value = 292
$ws.Range("K8").Value = 643
$ws.Range("K11").Value = 207
$ws.Range("K15").Value = 100
$ws.Range("K17").Value = 18
$ws.Range("K21").Value = 28
$ws.Range("K27").Value = 102
$ws.Range("K29").Value = 510
$ws.Range("K30").Value = 31
$ws.Range("K31").Value = 111
$ws.Range("K33").Value = 389
$ws.Range("K36").Value = 115
$ws.Range("K37").Value = 324
$ws.Range("K38").Value = 11
$ws.Range("K41").Value = 86
$ws.Range("K42").Value = 344
$ws.Range("K43").Value = 87
$ws.Range("K48").Value = 118
$ws.Range("K50").Value = 59
$ws.Range("K52").Value = 277
$ws.Range("K63").Value = 46
$ws.Range("K64").Value = 61
$ws.Range("K65").Value = 232
$ws.Range("K67").Value = 388
$ws.Range("K71").Value = 33
$ws.Range("K75").Value = 36
$ws.Range("K76").Value = 151
$ws.Range("K79").Value = 249
$ws.Range("K83").Value = 215
$ws.Range("K85").Value = 469
$ws.Range("K86").Value = 64
$ws.Range("K89").Value = 133
$ws.Range("K90").Value = 88
$ws.Range("K91").Value = 94
$ws.Range("K94").Value = 117
$ws.Range("K95").Value = 163
$ws.Range("K96").Value = 135
$ws.Range("K97").Value = 85
$ws.Range("K98").Value = 56
$ws.Range("K99").Value = 178
$ws.Range("K101").Value = 9841

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 100
$ws.Range("K3").Value = 88
$ws.Range("K5").Value = 13
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 292

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K3").Value = 55
$ws.Range("K7").Value = 207

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K3").Value = 43
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 133

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K3").Value = 162
$ws.Range("K6").Value = 106
$ws.Range("K7").Value = 469

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 75
$ws.Range("K3").Value = 69
$ws.Range("K6").Value = 113
$ws.Range("K7").Value = 277

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 187
$ws.Range("K7").Value = 643

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 80
$ws.Range("K3").Value = 67
$ws.Range("K7").Value = 215

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 141
$ws.Range("K6").Value = 111
$ws.Range("K7").Value = 389

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 51
$ws.Range("K3").Value = 57
$ws.Range("K7").Value = 163

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 111
$ws.Range("K6").Value = 102
$ws.Range("K7").Value = 324

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K3").Value = 54
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 232

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 178

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K2").Value = 9
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 126
$ws.Range("K6").Value = 112
$ws.Range("K7").Value = 388

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 139
$ws.Range("K3").Value = 169
$ws.Range("K4").Value = 29
$ws.Range("K6").Value = 161
$ws.Range("K7").Value = 510

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K3").Value = 23
$ws.Range("K7").Value = 118

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 90
$ws.Range("K7").Value = 151

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 87
$ws.Range("K3").Value = 110
$ws.Range("K7").Value = 344

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 42
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 94

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("K2").Value = 5
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 82
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 249

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 61

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K3").Value = 35
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 115

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 117

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K3").Value = 25
$ws.Range("K7").Value = 100

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K6").Value = 52
$ws.Range("K7").Value = 85

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K3").Value = 23
$ws.Range("K7").Value = 102

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K2").Value = 16
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 34
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K2").Value = 15
$ws.Range("K3").Value = 9
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range("K5").Value = 5
$ws.Range("K6").Value = 11
